# Add testing accuracy for bubble sheet:
# extend the question headers from Q1..Q16 to Q1..Q40 and append a second
# test-taker's ID + answer-key match results (1 = correct, 0 = incorrect).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- extend header row (1) with Q17 .. Q40, matching the existing header style ---
$headerStyleSource = $ws.Range("Q1")
for ($q = 17; $q -le 40; $q++) {
    $col = $q + 1   # Q1 is column 17, so Q17 -> column 18, etc.
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = "Q$q"
    $headerStyleSource.Copy()
    $cell.PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = $false

# --- update row 2 (the first respondent's ID + answers for Q1..Q16) ---
# Write the ID as text (not a number) without disturbing A2's existing
# (unstyled) cell format: build the text value in a scratch cell via a
# formula, then paste only its value into A2.
$scratch = $ws.Range("ZZ999")
$scratch.Formula = "=""4521"""
$scratch.Copy()
$ws.Range("A2").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()
$excel.CutCopyMode = $false

$row2Values = @(1,0,0,0,0,1,1,1,1,0,0,0,0,0,0,1,1,0,1,0,1,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

for ($i = 0; $i -lt $row2Values.Length; $i++) {
    $col = $i + 2   # B2 is column 2
    $ws.Cells.Item(2, $col).Value = $row2Values[$i]
}
